$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.024.95"
$ws.Range("E2").Value = "  +3.20%  "
$ws.Range("D3").Value = "3.097.03"
$ws.Range("E3").Value = "  +1.09%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +1.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.38"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.383"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.17%  "
$ws.Range("D12").Value = "3.631.21"
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.62%  "
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("D16").Value = "59.014.96"
$ws.Range("E16").Value = "  +2.81%  "
$ws.Range("D17").Value = "3.104.17"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "343.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.507"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").Value = "0.0₃0927"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.90%  "
$ws.Range("E29").Value = "  +2.55%  "
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("E31").Value = "  +2.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "155.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("E34").Value = "  +2.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.80%  "
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.75%  "
$ws.Range("D40").Value = "3.139.81"
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.666"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("E44").Value = "  +6.03%  "
$ws.Range("D45").Value = "2.295.79"
$ws.Range("E45").Value = "  +1.31%  "
$ws.Range("E46").Value = "  +1.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.76%  "
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.754"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "265.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +12.48%  "
